$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update % contribution figures for Design/Documentation/Implementation/Testing table (rows 4-7) ---
$ws.Range("E4").Value = 0.2
$ws.Range("F4").Value = 0.2

$ws.Range("E5").Value = 0.3
$ws.Range("F5").Value = 0.25

$ws.Range("E6").Value = 0.2
$ws.Range("F6").Value = 0.25

$ws.Range("E7").Value = 0.3
$ws.Range("F7").Value = 0.3
# Row 8 contains SUM formulas (E8/F8) that recalc automatically to 1

# --- Update file contribution hours for row 48 ---
$ws.Range("C48").Value = 8

# --- Add new Meetings Summary entry (row 49) ---
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)  # xlPasteFormats, copies the date number format (matches A38:A48 style)
$ws.Range("A49").Value = 45627
$ws.Range("B49").Value = "Developing "
$ws.Range("C49").Value = 3
$ws.Range("D49").Value = "Online Call"
$ws.Range("E49").Value = "Unavailable"
$ws.Range("F49").Value = "Unavailable"
$ws.Range("G49").Value = "Online Call"

# --- Update the active cell selection to reflect where the author last left off ---
$ws.Range("K45").Select()
